$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Reorder columns B..M to the new layout using whole-column cut/insert
#    (this preserves the exact bestFit width value that travels with each column)
$ops = @(
  @(8,2), @(9,3), @(12,4), @(13,5), @(7,6), @(8,7), @(9,8), @(10,9), @(11,10), @(12,11), @(13,12)
)
foreach ($op in $ops) {
    $src = $op[0]
    $dst = $op[1]
    $ws.Columns.Item($src).Cut()
    $ws.Columns.Item($dst).Insert()
}

# Column M (13) needs width "8" (recalculated bestFit after move) instead of the carried-over "8.28515625"
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667

# 2. Overwrite row 2 with the new "total score" values (exact, to avoid any precision loss from cut/paste)
$ws.Range("A2").Value = 42605.455625000002
$ws.Range("B2").Value = 13
$ws.Range("C2").Value = 87
$ws.Range("D2").Value = 97
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = "Random"

# 3. Delete row 3 (second, now-unneeded, Random sample)
$ws.Rows.Item(3).Delete()
